# Daily attendance processing - 2025-12-04 08:34:52
# Reorders the "Recorded By" (column G) contributor list for each session row so
# that the author names are listed in last-modified-first order. Cells whose
# "Recorded By" value only has a single contributor are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "" -and $val -like "*,*") {
        $parts = $val -split ","
        $n = $parts.Length
        $result = ""
        for ($i = $n - 1; $i -ge 0; $i--) {
            $piece = $parts[$i].Trim()
            if ($i -eq ($n - 1)) {
                $result = $piece
            } else {
                $result = $result + ", " + $piece
            }
        }
        $cell.Value2 = $result
    }
}
